$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D, shifting existing D:K data to E:L
$ws.Range("D:D").Insert()

# Copy formatting from column E (the shifted original column D) into new column D
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with the new period (2018-12-31) figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1279300
$ws.Range("D9").Value = 1164600
$ws.Range("D10").Value = 114700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1231000
$ws.Range("D18").Value = 48300
$ws.Range("D20").Value = 100
$ws.Range("D21").Value = 52200
$ws.Range("D22").Value = 5200
$ws.Range("D23").Value = 43300
$ws.Range("D24").Value = 10000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 33300
$ws.Range("D27").Value = 33300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 2000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -100
$ws.Range("D33").Value = 35300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 35300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 66100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 232100
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 15100
$ws.Range("D46").Value = 313300
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 13400
$ws.Range("D49").Value = 242200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 3200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 572200
$ws.Range("D57").Value = 156400
$ws.Range("D58").Value = 4500
$ws.Range("D59").Value = 64100
$ws.Range("D60").Value = 225000
$ws.Range("D61").Value = 69100
$ws.Range("D62").Value = 56800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 350900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 152600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 221300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 35300
$ws.Range("D83").Value = 3800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 40100
$ws.Range("D91").Value = -10000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -46800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -3300
$ws.Range("D101").Value = -1200
$ws.Range("D102").Value = -11300
